$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update period headers (row 8) and publish dates (row 9) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-02-11 (8)"
$ws.Range("E9").Value = "1400-02-08 (8)"
$ws.Range("F9").Value = "1401-02-06 (9)"
$ws.Range("G9").Value = "1402-02-10 (8)"
$ws.Range("H9").Value = "1402-02-10 (2)"

# --- Update financial data rows (shift columns left, append new period column H) ---
$ws.Range("D11").Value = 16983
$ws.Range("E11").Value = 20591
$ws.Range("F11").Value = 15077
$ws.Range("G11").Value = 31495
$ws.Range("H11").Value = 31261

$ws.Range("D12").Value = -13547
$ws.Range("E12").Value = -16454
$ws.Range("F12").Value = -10988
$ws.Range("G12").Value = -24925
$ws.Range("H12").Value = -24885

$ws.Range("D13").Value = "-"
$ws.Range("E13").Value = 4137
$ws.Range("F13").Value = 4089
$ws.Range("G13").Value = 6570
$ws.Range("H13").Value = 6375

$ws.Range("D14").Value = -554
$ws.Range("E14").Value = -530
$ws.Range("F14").Value = -518
$ws.Range("G14").Value = -621
$ws.Range("H14").Value = -612

$ws.Range("D15").Value = -172
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = "-"

$ws.Range("D16").Value = -55
$ws.Range("E16").Value = 42
$ws.Range("F16").Value = 135
$ws.Range("G16").Value = 438
$ws.Range("H16").Value = 94

$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = 3649
$ws.Range("F17").Value = 3706
$ws.Range("G17").Value = 6387
$ws.Range("H17").Value = 5858

$ws.Range("D18").Value = -410
$ws.Range("E18").Value = -129
$ws.Range("F18").Value = -40
$ws.Range("G18").Value = -54
$ws.Range("H18").Value = -38

$ws.Range("D19").Value = 645
$ws.Range("E19").Value = -3
$ws.Range("F19").Value = -195
$ws.Range("G19").Value = 65
$ws.Range("H19").Value = -36

$ws.Range("D20").Value = "-"
$ws.Range("E20").Value = 3517
$ws.Range("F20").Value = 3471
$ws.Range("G20").Value = 6398
$ws.Range("H20").Value = 5783

$ws.Range("D21").Value = -692
$ws.Range("E21").Value = -788
$ws.Range("F21").Value = -503
$ws.Range("G21").Value = -1071
$ws.Range("H21").Value = -594

$ws.Range("D22").Value = "-"
$ws.Range("E22").Value = 2729
$ws.Range("F22").Value = 2968
$ws.Range("G22").Value = 5327
$ws.Range("H22").Value = 5189

$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = 2729
$ws.Range("F24").Value = 2968
$ws.Range("G24").Value = 5327
$ws.Range("H24").Value = 5189

$ws.Range("D25").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0

$ws.Range("D26").Value = 4691
$ws.Range("E26").Value = 3699
$ws.Range("F26").Value = 2099
$ws.Range("G26").Value = 1798
$ws.Range("H26").Value = 2692

$ws.Range("D27").Value = "-"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
